# Ablation-table refresh: the MAE/MAPE columns for the "variant" rows (9-13)
# on Sheet1 were regenerated (charts now cache local values instead of an
# external link), so the underlying cells need the new numbers. The source
# layout being pasted in also carries two blank helper columns after each
# pair of value columns (mirrors the header block in rows 1-6), so those
# land as genuinely blank-but-present cells - reproduce that with a
# copy/paste from a scratch area instead of a plain value assignment
# (a plain ClearContents / $null assign removes the cell entirely instead
# of leaving an empty-but-present cell behind).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Refresh the three 5-row data blocks (rows 9-13) that feed chart2-4.
#    Each block is rotated up by one row (row 9 gets what used to be row
#    10's numbers, ..., row 13 wraps around to the old row 9 numbers) and
#    gains two blank spacer columns, matching the pasted-in source range.
# ---------------------------------------------------------------------

$scratch = 200

# Block 1: G:J (MAE/MAPE pair 2) -----------------------------------------
$ws.Range("G$scratch").Value  = 17.86
$ws.Range("H$scratch").Value  = 17.25
$ws.Range("G$($scratch+1)").Value = 18.42
$ws.Range("H$($scratch+1)").Value = 17.91
$ws.Range("G$($scratch+2)").Value = 17.15
$ws.Range("H$($scratch+2)").Value = 16.58
$ws.Range("G$($scratch+3)").Value = 16.37
$ws.Range("H$($scratch+3)").Value = 15.74
$ws.Range("G$($scratch+4)").Value = 15.23
$ws.Range("H$($scratch+4)").Value = 14.83
# I/J columns stay untouched (blank) in the scratch rows on purpose.
$ws.Range("G$scratch`:J$($scratch+4)").Copy($ws.Range("G9:J13"))
$ws.Range("G$scratch`:J$($scratch+4)").Clear()

# Block 2: L:O (MAE/MAPE pair 3) ------------------------------------------
$ws.Range("L$scratch").Value  = 21.59
$ws.Range("M$scratch").Value  = 14.62
$ws.Range("L$($scratch+1)").Value = 22.17
$ws.Range("M$($scratch+1)").Value = 15.13
$ws.Range("L$($scratch+2)").Value = 20.94
$ws.Range("M$($scratch+2)").Value = 13.87
$ws.Range("L$($scratch+3)").Value = 19.78
$ws.Range("M$($scratch+3)").Value = 13.15
$ws.Range("L$($scratch+4)").Value = 18.65
$ws.Range("M$($scratch+4)").Value = 12.39
# N/O columns stay untouched (blank) in the scratch rows on purpose.
$ws.Range("L$scratch`:O$($scratch+4)").Copy($ws.Range("L9:O13"))
$ws.Range("L$scratch`:O$($scratch+4)").Clear()

# Block 3: Q:R (MAE/MAPE pair 4, no spacer columns after it) --------------
$ws.Range("Q$scratch").Value  = 17.34
$ws.Range("R$scratch").Value  = 11.28
$ws.Range("Q$($scratch+1)").Value = 17.96
$ws.Range("R$($scratch+1)").Value = 11.75
$ws.Range("Q$($scratch+2)").Value = 16.82
$ws.Range("R$($scratch+2)").Value = 10.63
$ws.Range("Q$($scratch+3)").Value = 15.69
$ws.Range("R$($scratch+3)").Value = 10.12
$ws.Range("Q$($scratch+4)").Value = 14.71
$ws.Range("R$($scratch+4)").Value = 9.45
$ws.Range("Q$scratch`:R$($scratch+4)").Copy($ws.Range("Q9:R13"))
$ws.Range("Q$scratch`:R$($scratch+4)").Clear()

# ---------------------------------------------------------------------
# 2) Restore VBA codeNames (workbook + each sheet) so the project keeps
#    its stable module identity instead of the blank default.
# ---------------------------------------------------------------------
$wb.CodeName = "ThisWorkbook"
$ws.CodeName = "Sheet1"
$wb.Worksheets.Item("Sheet2").CodeName = "Sheet2"
$wb.Worksheets.Item("Sheet3").CodeName = "Sheet3"

# ---------------------------------------------------------------------
# 3) Move the viewport/selection on Sheet1 back into view of the table
#    (scrolled up a few rows, selection moved from O17 to O19).
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.TopLeftCell = $ws.Range("E1")
$ws.Range("O19").Select()
